{"js": "// Replace each old three-digit-division answer with its new value.\n// The mapping is 1:1 and every \"old\" string is unique in the document,\n// so a body.search() + insertText(..., \"Replace\") pair is unambiguous\n// for each entry.\nconst pairs = [\n  [\"993\u00f74=248, 1\", \"604\u00f77=86, 2\"],\n  [\"899\u00f76=149, 5\", \"357\u00f78=44, 5\"],\n  [\"309\u00f76=51, 3\", \"407\u00f73=135, 2\"],\n  [\"635\u00f72=317, 1\", \"149\u00f75=29, 4\"],\n  [\"801\u00f76=133, 3\", \"231\u00f78=28, 7\"],\n  [\"993\u00f77=141, 6\", \"624\u00f78=78, 0\"],\n  [\"362\u00f74=90, 2\", \"538\u00f75=107, 3\"],\n  [\"911\u00f73=303, 2\", \"675\u00f73=225, 0\"],\n  [\"391\u00f76=65, 1\", \"810\u00f78=101, 2\"],\n  [\"437\u00f74=109, 1\", \"714\u00f78=89, 2\"],\n  [\"988\u00f78=123, 4\", \"721\u00f78=90, 1\"],\n  [\"470\u00f76=78, 2\", \"838\u00f73=279, 1\"],\n  [\"249\u00f72=124, 1\", \"884\u00f78=110, 4\"],\n  [\"900\u00f76=150, 0\", \"906\u00f77=129, 3\"],\n  [\"167\u00f78=20, 7\", \"746\u00f74=186, 2\"],\n  [\"125\u00f77=17, 6\", \"390\u00f79=43, 3\"],\n  [\"526\u00f74=131, 2\", \"467\u00f77=66, 5\"],\n  [\"967\u00f78=120, 7\", \"407\u00f72=203, 1\"],\n  [\"590\u00f79=65, 5\", \"112\u00f75=22, 2\"],\n  [\"291\u00f78=36, 3\", \"130\u00f77=18, 4\"],\n  [\"667\u00f73=222, 1\", \"929\u00f74=232, 1\"],\n  [\"706\u00f77=100, 6\", \"949\u00f79=105, 4\"],\n  [\"835\u00f75=167, 0\", \"945\u00f73=315, 0\"],\n  [\"645\u00f79=71, 6\", \"988\u00f77=141, 1\"],\n  [\"660\u00f72=330, 0\", \"235\u00f78=29, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old three-digit-division answer with its new value.\n# The mapping is 1:1 and every \"old\" string is unique in the document,\n# so a plain Find/Replace (no wildcards) for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('993\u00f74=248, 1', '604\u00f77=86, 2'),\n    @('899\u00f76=149, 5', '357\u00f78=44, 5'),\n    @('309\u00f76=51, 3', '407\u00f73=135, 2'),\n    @('635\u00f72=317, 1', '149\u00f75=29, 4'),\n    @('801\u00f76=133, 3', '231\u00f78=28, 7'),\n    @('993\u00f77=141, 6', '624\u00f78=78, 0'),\n    @('362\u00f74=90, 2', '538\u00f75=107, 3'),\n    @('911\u00f73=303, 2', '675\u00f73=225, 0'),\n    @('391\u00f76=65, 1', '810\u00f78=101, 2'),\n    @('437\u00f74=109, 1', '714\u00f78=89, 2'),\n    @('988\u00f78=123, 4', '721\u00f78=90, 1'),\n    @('470\u00f76=78, 2', '838\u00f73=279, 1'),\n    @('249\u00f72=124, 1', '884\u00f78=110, 4'),\n    @('900\u00f76=150, 0', '906\u00f77=129, 3'),\n    @('167\u00f78=20, 7', '746\u00f74=186, 2'),\n    @('125\u00f77=17, 6', '390\u00f79=43, 3'),\n    @('526\u00f74=131, 2', '467\u00f77=66, 5'),\n    @('967\u00f78=120, 7', '407\u00f72=203, 1'),\n    @('590\u00f79=65, 5', '112\u00f75=22, 2'),\n    @('291\u00f78=36, 3', '130\u00f77=18, 4'),\n    @('667\u00f73=222, 1', '929\u00f74=232, 1'),\n    @('706\u00f77=100, 6', '949\u00f79=105, 4'),\n    @('835\u00f75=167, 0', '945\u00f73=315, 0'),\n    @('645\u00f79=71, 6', '988\u00f77=141, 1'),\n    @('660\u00f72=330, 0', '235\u00f78=29, 3'),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n\n"}
